$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF"
# Copy formatting from the existing header cell (H1) so the new header
# cells share the same style (bold, centered, bordered) as the others.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-30
$values = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(8, 8)
    5  = @(9, 9)
    6  = @(8, 8)
    7  = @(9, 9)
    8  = @(9, 9)
    9  = @(8, 9)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(7, 8)
    14 = @(9, 9)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(6, 6)
    19 = @(8, 8)
    20 = @(8, 8)
    21 = @(9, 9)
    22 = @(6, 6)
    23 = @(7, 7)
    24 = @(9, 9)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(5, 5)
    29 = @(9, 9)
    30 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
